$wb = $excel.ActiveWorkbook

# --- Sheet 1 (展览) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 7630
$ws1.Range("F5").Value = 931
$ws1.Range("F7").Value = 796
$ws1.Range("F9").Value = 88
$ws1.Range("F13").Value = 3066
$ws1.Range("F14").Value = 191
$ws1.Range("F16").Value = 718
$ws1.Range("F17").Value = 750
$ws1.Range("F19").Value = 449
$ws1.Range("F21").Value = 211
$ws1.Range("F22").Value = 212
$ws1.Range("F23").Value = 264
$ws1.Range("F26").Value = 94
$ws1.Range("F27").Value = 248
$ws1.Range("F30").Value = 492
$ws1.Range("F31").Value = 451
$ws1.Range("F32").Value = 27
$ws1.Range("F35").Value = 83

# --- Sheet 2 (演出) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 35
$ws2.Range("G2").Value = 108

# --- Sheet 4 (全部类型) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 35
$ws4.Range("G3").Value = 108
$ws4.Range("F6").Value = 7630
$ws4.Range("F8").Value = 931
$ws4.Range("F10").Value = 796
$ws4.Range("F12").Value = 88
$ws4.Range("F17").Value = 3066
$ws4.Range("F18").Value = 191
$ws4.Range("F21").Value = 718
$ws4.Range("F22").Value = 750
$ws4.Range("F25").Value = 449
$ws4.Range("F27").Value = 211
$ws4.Range("F28").Value = 212
$ws4.Range("F29").Value = 264
$ws4.Range("F32").Value = 94
$ws4.Range("F33").Value = 248
$ws4.Range("F36").Value = 492
$ws4.Range("F37").Value = 451
$ws4.Range("F38").Value = 27
$ws4.Range("F41").Value = 83
